# Auto-generated Word COM-interop PowerShell script
$d = $word.ActiveDocument

# 1) Update the title/date paragraph
$d.Content.Find.Execute("2023-11-07 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-08 Wednesday", 2) | Out-Null

# 2) Update each arithmetic expression cell in the table, in row-major order
$t = $d.Tables.Item(1)
$newValues = @(
    "56+9=",
    "30-24=",
    "39+56=",
    "61-18=",
    "7+36=",
    "49+4=",
    "32-16=",
    "69+28=",
    "39+22=",
    "74-66=",
    "47-19=",
    "34+47=",
    "16+26=",
    "52-49=",
    "33-15=",
    "94-39=",
    "17+56=",
    "44+39=",
    "23+69=",
    "4+39=",
    "37+58=",
    "70-41=",
    "48+18=",
    "6+66=",
    "51-25=",
    "48+24=",
    "71-59=",
    "12-4=",
    "78-69=",
    "7+68=",
    "63-24=",
    "91-27=",
    "87+7=",
    "58+14=",
    "25-17=",
    "26+35=",
    "30-12=",
    "50-13=",
    "45+8=",
    "44+19=",
    "24+59=",
    "86-68=",
    "7+84=",
    "96-57=",
    "8+67=",
    "26+56=",
    "3+8=",
    "50-33=",
    "39+59=",
    "15+18=",
    "87-58=",
    "28+27=",
    "33+18=",
    "88+5=",
    "94-59=",
    "50-25=",
    "76-27=",
    "50-9=",
    "47+34=",
    "69+8=",
    "50-28=",
    "92-55=",
    "51-46=",
    "53+28=",
    "38+18=",
    "80-63=",
    "40-9=",
    "83-5=",
    "74-9=",
    "8+4=",
    "44-37=",
    "68+8=",
    "44+37=",
    "52-29=",
    "44-9=",
    "26+38=",
    "35+36=",
    "30-15=",
    "74+7=",
    "4+39=",
    "72-45=",
    "25+9=",
    "33-7=",
    "24+49=",
    "27+27=",
    "26+25=",
    "27+46=",
    "26+5=",
    "9+36=",
    "87-68=",
    "62-37=",
    "80-5=",
    "22+69=",
    "55-17=",
    "55+28=",
    "4+17=",
    "51-35=",
    "5+47=",
    "25+66=",
    "3+9="
)

$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Done. Updated $idx cells."
